# Updated capital structure database
# Refreshes the financial metrics for the two Iceland / Air Transport rows
# (row 2 = industry aggregate, row 3 = Icelandair Group hf.) with newly
# recomputed values, and drops the now-unused buybacks_cash_returned (T)
# column for both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "D"  = -0.09080000000000001
    "G"  = 0.1539350180505415
    "H"  = 0.1539350180505415
    "I"  = -0.3211552346570397
    "J"  = -0.3211552346570397
    "K"  = -313.5
    "L"  = -0.4527075812274368
    "M"  = 0
    "N"  = 0
    "O"  = 0
    "S"  = 0
    "U"  = 191
    "V"  = 0.5210038188761593
    "W"  = -0.6256236280183596
    "X"  = 0.1058043949302045
    "Y"  = -0.7314280229485641
    "Z"  = 0.9041650345998171
    "AA" = -0.2903773338555947
    "AB" = 0.06394524569979672
    "AC" = -0.3543225795553914
    "AD" = 411.3
    "AE" = 0
    "AF" = 411.3
    "AG" = 220.3
    "AH" = 0.5287311993829541
    "AI" = 0.5839011925042589
    "AJ" = 0.3753620719032203
    "AK" = 0.4291001168679392
    "AL" = 35.2
    "AM" = 32.19
    "AN" = -6.366873065015481
    "AO" = -6.318181818181817
    "AP" = -3.410216718266254
    "AQ" = -6.908977943460701
}

foreach ($row in 2, 3) {
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
    # buybacks_cash_returned column is no longer reported - remove it
    $ws.Range("T$row").ClearContents()
}
